$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 corresponds to year 2025: update total_customers, new_customers,
# new_rate and returning_rate to reflect the updated counts.
$ws.Range("C6").Value = 408
$ws.Range("E6").Value = 101
$ws.Range("G6").Value = 24.75490196078432
$ws.Range("H6").Value = 75.24509803921569
